$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be
# auto-converted to numbers (losing exact display text, trailing zeros, etc.)
$textCells = @("D5","D6","D8","D12","D14","D16","D20","D21","D23","D24","D27","D29","D30","D31","D32","D35","D39","D40","D41","D43","D45","D46","D49","D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "67.366.14"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "2.636.70"
$ws.Range("E3").Value = "  -2.68%  "
$ws.Range("D5").Value = "596.16"
$ws.Range("E5").Value = "  -2.24%  "
$ws.Range("D6").Value = "168.38"
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.534"
$ws.Range("E8").Value = "  -3.62%  "
$ws.Range("D9").Value = "2.636.19"
$ws.Range("E9").Value = "  -2.68%  "
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("D12").Value = "0.362"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").Value = "27.88"
$ws.Range("E14").Value = "  -1.73%  "
$ws.Range("D15").Value = "3.114.22"
$ws.Range("E15").Value = "  -2.80%  "
$ws.Range("D16").Value = "0.0000182"
$ws.Range("E16").Value = "  -2.63%  "
$ws.Range("D17").Value = "67.120.00"
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("D18").Value = "2.627.62"
$ws.Range("E18").Value = "  -2.63%  "
$ws.Range("E19").Value = "  +3.57%  "
$ws.Range("D20").Value = "8.12"
$ws.Range("E20").Value = "  +6.54%  "
$ws.Range("D21").Value = "360.39"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("E22").Value = "  -2.66%  "
$ws.Range("D23").Value = "4.70"
$ws.Range("E23").Value = "  -4.22%  "
$ws.Range("D24").Value = "10.90"
$ws.Range("E24").Value = "  +9.06%  "
$ws.Range("E25").Value = "  -5.17%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "70.54"
$ws.Range("E27").Value = "  -3.33%  "
$ws.Range("D28").Value = "2.769.35"
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0000102"
$ws.Range("E30").Value = "  -2.44%  "
$ws.Range("D31").Value = "556.88"
$ws.Range("E31").Value = "  -3.53%  "
$ws.Range("D32").Value = "7.95"
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("E33").Value = "  -2.91%  "
$ws.Range("E34").Value = "  -3.41%  "
$ws.Range("D35").Value = "0.137"
$ws.Range("E35").Value = "  +4.36%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  -4.52%  "
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("D39").Value = "19.22"
$ws.Range("E39").Value = "  -3.10%  "
$ws.Range("D40").Value = "0.367"
$ws.Range("E40").Value = "  -2.74%  "
$ws.Range("D41").Value = "5.22"
$ws.Range("E41").Value = "  -2.97%  "
$ws.Range("E42").Value = "  -2.99%  "
$ws.Range("D43").Value = "17.94"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "2.48"
$ws.Range("E45").Value = "  -4.43%  "
$ws.Range("D46").Value = "40.20"
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("E47").Value = "  -2.79%  "
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("D49").Value = "152.51"
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("D50").Value = "3.84"
$ws.Range("E50").Value = "  -1.13%  "
